# Regenerate sval data to filter save games.
# Updates the numeric values in columns B-E and G for rows 2-7
# while leaving the date strings in column A and the Win flags
# in column F untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 1.455362044514542;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 3.754798637575387 }
    3 = @{ B = 3.286832544864788;  C = 0.306821227259698;   D = 3.537761648806719;  E = 0.4942365360607697; G = 7.625651956991975 }
    4 = @{ B = 1.455362044514542;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 3.754798637575387 }
    5 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    6 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    7 = @{ B = 0.6606524410359556; C = 1.655778082260271;  D = 261.3203778131603;  E = 1133.036916526867;  G = 1396.673724863324 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
